# Commit: "Reduced Order history and category History size in ui.js"
#
# The "Order History" sheet is regenerated from a smaller page size: the
# two existing data rows (bill #432/#433, placed on 45903.x) are replaced
# by eight rows for bill #586..#593 (placed on 45692.x), growing the used
# range from A1:I3 to A1:I9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bill rows, in the order they should land in rows 2..9.
$rows = @(
    @(593, 45692.00011574074, "Ajay Francis Anchan", 8, 60,  0.9,  0.9,  1.8,   "Chicken Burger (x1)"),
    @(592, 45692.00011574074, "Ajay Francis Anchan", 7, 200, 3.5,  3.5,  7,     "Chicken Cheese Burger (x1), Vanilla Shake (x2), Chicken Wrap (x1)"),
    @(591, 45692.00011574074, "Ajay Francis Anchan", 6, 80,  0.8,  0.8,  1.6,   "Belgian Coffee (x2)"),
    @(590, 45692.00011574074, "Ajay Francis Anchan", 5, 140, 2.2,  2.2,  2,     "Peri Peri Fries (x2), Watermelon Juice (x1), Belgian Coffee (x1), Lime Juice (x1)"),
    @(589, 45692.00011574074, "Ajay Francis Anchan", 4, 220, 3.1,  3.1,  6.2,   "Chicken Burger (x2), Mango Lassi (x1), Strawberry Lassi (x1)"),
    @(588, 45692.00011574074, "Ajay Francis Anchan", 3, 180, 1.2,  1.2,  4.2,   "Vanilla Shake (x3), Chicken Cheese Burger (x1)"),
    @(587, 45692.00011574074, "Ajay Francis Anchan", 2, 370, 5.03, 5.03, 10.05, "Oreo Shake (x1), Vanilla Shake (x1), Chicken Wrap (x2), Butterscotch Lassi (x2), Strawberry Lassi (x1)"),
    @(586, 45692.00011574074, "Ajay Francis Anchan", 1, 300, 3.6,  3.6,  9,     "Chicken Burger (x2), Chicken Cheese Burger (x1), Butterscotch Lassi (x2)")
)

$r = 2
foreach ($row in $rows) {
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col - 1]
    }
    $r++
}

# B2 already carries the date number format (style index 1). Propagate
# that same format down the freshly written date column B3:B9 so every
# row matches (rather than leaving the new cells as "General").
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3:B9").PasteSpecial(-4122) | Out-Null
